$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Cuoi ky"
$ws.Range("C1").Value = "Giua ky"
$ws.Range("D1").Value = "Total"

# Row 2
$ws.Range("A2").Value = "Nguyen Van A"
$ws.Range("B2").Value = "'20"
$ws.Range("C2").Value = "'20"
$ws.Range("D2").Value = "'40"

# Row 3
$ws.Range("A3").Value = "Nguyen Van B"
$ws.Range("B3").Value = "'30"
$ws.Range("C3").Value = "'30"
$ws.Range("D3").Value = "'60"

# Row 4 (new row)
$ws.Range("A4").Value = "Student"
$ws.Range("B4").Value = "'50"
$ws.Range("C4").Value = "'50"
$ws.Range("D4").Value = "'100"
